# Update leve-crafting profit figures (H/I/J/K/L/M/N columns) across
# several sheets with refreshed market-board pricing data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 833.6
$ws.Range("I28").Value = 822.61536
$ws.Range("K28").Value = 822.61536
$ws.Range("M28").Value = -337.61536

# Row 41
$ws.Range("H41").Value = 152.75
$ws.Range("I41").Value = 56.2
$ws.Range("J41").Value = 313.66666
$ws.Range("K41").Value = 56.2
$ws.Range("L41").Value = 313.66666
$ws.Range("M41").Value = 383.8
$ws.Range("N41").Value = -1193.66666

# Row 46
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").Value = $null

# Row 51
$ws.Range("H51").Value = 11562
$ws.Range("J51").Value = 14333.333
$ws.Range("L51").Value = 14333.333
$ws.Range("N51").Value = -15301.333

# Row 60
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").Value = $null

# Row 82
$ws.Range("H82").Value = 1533.3334
$ws.Range("I82").Value = 1533.3334
$ws.Range("K82").Value = 4600.0002
$ws.Range("M82").Value = -4194.0002

# Row 85
$ws.Range("H85").Value = 1533.3334
$ws.Range("I85").Value = 1533.3334
$ws.Range("K85").Value = 4600.0002
$ws.Range("M85").Value = -3196.0002

# Row 96
$ws.Range("H96").Value = 2094.3333
$ws.Range("I96").Value = 2141.5
$ws.Range("J96").Value = 2000
$ws.Range("K96").Value = 6424.5
$ws.Range("L96").Value = 6000
$ws.Range("M96").Value = -5051.5
$ws.Range("N96").Value = -8746

# Row 100
$ws.Range("H100").Value = 1657.5714
$ws.Range("I100").Value = 1698.4
$ws.Range("J100").Value = 1555.5
$ws.Range("K100").Value = 1698.4
$ws.Range("L100").Value = 1555.5
$ws.Range("M100").Value = -1157.4
$ws.Range("N100").Value = -2637.5

# Row 113
$ws.Range("H113").Value = 7930.1177
$ws.Range("I113").Value = 10786.714
$ws.Range("J113").Value = 5930.5
$ws.Range("K113").Value = 10786.714
$ws.Range("L113").Value = 5930.5
$ws.Range("M113").Value = -7532.714
$ws.Range("N113").Value = -12438.5

# Row 138
$ws.Range("H138").Value = 1814.0769
$ws.Range("I138").Value = 1308.3
$ws.Range("K138").Value = 3924.9
$ws.Range("M138").Value = 1215.1

$ws = $wb.Worksheets.Item("ARM")
# Row 122
$ws.Range("H122").Value = 1799
$ws.Range("I122").Value = 1799
$ws.Range("K122").Value = 5397
$ws.Range("M122").Value = -2947

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 656.7273
$ws.Range("J20").Value = 869.4
$ws.Range("L20").Value = 869.4
$ws.Range("N20").Value = -1363.4

# Row 107
$ws.Range("H107").Value = 997.25
$ws.Range("I107").Value = 997.25
$ws.Range("K107").Value = 997.25
$ws.Range("M107").Value = 922.75

# Row 134
$ws.Range("H134").Value = 537
$ws.Range("I134").Value = 537
$ws.Range("K134").Value = 1611
$ws.Range("M134").Value = 924

$ws = $wb.Worksheets.Item("CRP")
# Row 23
$ws.Range("H23").Value = 75000
$ws.Range("I23").Value = 75000
$ws.Range("J23").Value = 75000
$ws.Range("K23").Value = 75000
$ws.Range("L23").Value = 75000
$ws.Range("M23").Value = -74760
$ws.Range("N23").Value = -75480

# Row 27
$ws.Range("H27").Value = 75000
$ws.Range("I27").Value = 75000
$ws.Range("J27").Value = 75000
$ws.Range("K27").Value = 75000
$ws.Range("L27").Value = 75000
$ws.Range("M27").Value = -74808
$ws.Range("N27").Value = -75384

# Row 55
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").Value = $null

# Row 99
$ws.Range("H99").Value = 3252
$ws.Range("I99").Value = 2628
$ws.Range("J99").Value = 4500
$ws.Range("K99").Value = 2628
$ws.Range("L99").Value = 4500
$ws.Range("M99").Value = -1130
$ws.Range("N99").Value = -7496

# Row 126
$ws.Range("H126").Value = 3252
$ws.Range("I126").Value = 2628
$ws.Range("J126").Value = 4500
$ws.Range("K126").Value = 7884
$ws.Range("L126").Value = 13500
$ws.Range("M126").Value = -5414
$ws.Range("N126").Value = -18440

# Row 134
$ws.Range("H134").Value = 2355.5
$ws.Range("I134").Value = 2355.5
$ws.Range("K134").Value = 7066.5
$ws.Range("M134").Value = -4531.5

$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 420
$ws.Range("J34").Value = 600
$ws.Range("L34").Value = 1800
$ws.Range("N34").Value = -1968

# Row 130
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = $null

# Row 132
$ws.Range("H132").Value = 3383.3333
$ws.Range("J132").Value = 3500
$ws.Range("L132").Value = 31500
$ws.Range("N132").Value = -36560

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 122.61539
$ws.Range("I2").Value = 60.333332
$ws.Range("J2").Value = 141.3
$ws.Range("K2").Value = 60.333332
$ws.Range("L2").Value = 141.3
$ws.Range("M2").Value = 52.666668
$ws.Range("N2").Value = -367.3

# Row 43
$ws.Range("H43").Value = 15004
$ws.Range("I43").Value = 20008
$ws.Range("J43").Value = 10000
$ws.Range("K43").Value = 20008
$ws.Range("L43").Value = 10000
$ws.Range("M43").Value = -19857
$ws.Range("N43").Value = -10302

# Row 46
$ws.Range("H46").Value = 9999.799999999999
$ws.Range("I46").Value = 9999
$ws.Range("K46").Value = 9999
$ws.Range("M46").Value = -9843

# Row 57
$ws.Range("H57").Value = 17500
$ws.Range("I57").Value = 10000
$ws.Range("K57").Value = 10000
$ws.Range("M57").Value = -9180

# Row 97
$ws.Range("H97").Value = 2783.3333
$ws.Range("I97").Value = 2776.2778
$ws.Range("J97").Value = 2804.5
$ws.Range("K97").Value = 2776.2778
$ws.Range("L97").Value = 2804.5
$ws.Range("M97").Value = -2280.2778
$ws.Range("N97").Value = -3796.5

# Row 113
$ws.Range("H113").Value = 4351.8
$ws.Range("I113").Value = 2586.6667
$ws.Range("J113").Value = 6999.5
$ws.Range("K113").Value = 2586.6667
$ws.Range("L113").Value = 6999.5
$ws.Range("M113").Value = -416.6667000000002
$ws.Range("N113").Value = -11339.5

# Row 122
$ws.Range("H122").Value = 1403.8334
$ws.Range("I122").Value = 1403.8334
$ws.Range("K122").Value = 4211.5002
$ws.Range("M122").Value = -1761.5002

# Row 132
$ws.Range("H132").Value = 3022.4614
$ws.Range("I132").Value = 3022.4614
$ws.Range("K132").Value = 9067.3842
$ws.Range("M132").Value = -6537.3842

# Row 134
$ws.Range("H134").Value = 43442
$ws.Range("J134").Value = 43442
$ws.Range("L134").Value = 130326
$ws.Range("N134").Value = -135396

$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Range("H136").Value = 4278.2856
$ws.Range("I136").Value = 3737
$ws.Range("K136").Value = 11211
$ws.Range("M136").Value = -8661

$ws = $wb.Worksheets.Item("WVR")
# Row 26
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").Value = $null

# Row 29
$ws.Range("H29").Value = 500
$ws.Range("I29").Value = 500
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 500
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -210
$ws.Range("N29").Value = $null
